$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row (Chinese field labels -> canonical English field names) ---
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "area"
$ws.Range("D1").Value = "share_portion"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "acquire_value"

# --- New metadata header columns (I..O) ---
$ws.Range("I1").Value = "property_category"
$ws.Range("J1").Value = "category"
$ws.Range("K1").Value = "date"
$ws.Range("L1").Value = "legislator_name"
$ws.Range("M1").Value = "legislator_id"
$ws.Range("N1").Value = "source_file"
$ws.Range("O1").Value = "index"

# Give the new header cells (I1:O1) the same look as the existing header cells
# (bold, centered, top-aligned, thin border) so they match B1:H1.
$newHeaders = $ws.Range("I1:O1")
$newHeaders.Borders.LineStyle = 1
$newHeaders.Font.Bold = $true
$newHeaders.HorizontalAlignment = -4108
$newHeaders.VerticalAlignment = -4160

# --- Clean up land-plot names: drop the leading star/dash and the stray space before 地號 ---
$ws.Range("B2").Value = "臺南市東區德高段04880023地號"
$ws.Range("B3").Value = "臺南市東區德高段04880090地號"

# --- Clean up the registration-date strings: drop the stray internal space ---
$ws.Range("F2").Value = "99年12月13日"
$ws.Range("F3").Value = "99年12月13曰"

# --- Populate the new metadata columns for each data row ---
$ws.Range("I2").Value = "land"
$ws.Range("J2").Value = "normal"
$ws.Range("K2").Value = "'2010-12-24"
$ws.Range("K2").ClearFormats()
$ws.Range("L2").Value = "陳淑慧"
$ws.Range("M2").Value = 1720
$ws.Range("N2").Value = "tmpdd71"
$ws.Range("O2").Value = 14

$ws.Range("I3").Value = "land"
$ws.Range("J3").Value = "normal"
$ws.Range("K3").Value = "'2010-12-24"
$ws.Range("K3").ClearFormats()
$ws.Range("L3").Value = "陳淑慧"
$ws.Range("M3").Value = 1720
$ws.Range("N3").Value = "tmpdd71"
$ws.Range("O3").Value = 15
